$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/4/2023  Through  12/10/2023"

# --- Weekly crime-stat grid updates (rows 14-30, cols C-N) ---
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 4).NumberFormat = "#,##0"
$ws.Cells.Item(14, 5).Value = -100
$ws.Cells.Item(14, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(14, 7).Value = 2
$ws.Cells.Item(14, 7).NumberFormat = "#,##0"
$ws.Cells.Item(14, 10).Value = 7
$ws.Cells.Item(14, 10).NumberFormat = "#,##0"
$ws.Cells.Item(14, 11).Value = 71.428571428571
$ws.Cells.Item(14, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 6).NumberFormat = "#,##0"
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 7).NumberFormat = "#,##0"
$ws.Cells.Item(15, 8).Value = -66.666666666666
$ws.Cells.Item(15, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(15, 10).Value = 40
$ws.Cells.Item(15, 10).NumberFormat = "#,##0"
$ws.Cells.Item(15, 11).Value = -22.5
$ws.Cells.Item(15, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(15, 13).Value = 47.619047619047
$ws.Cells.Item(15, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(15, 14).Value = -57.534246575342
$ws.Cells.Item(15, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(16, 3).Value = 7
$ws.Cells.Item(16, 3).NumberFormat = "#,##0"
$ws.Cells.Item(16, 4).Value = 9
$ws.Cells.Item(16, 4).NumberFormat = "#,##0"
$ws.Cells.Item(16, 5).Value = -22.222222222222
$ws.Cells.Item(16, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(16, 6).Value = 32
$ws.Cells.Item(16, 6).NumberFormat = "#,##0"
$ws.Cells.Item(16, 7).Value = 33
$ws.Cells.Item(16, 7).NumberFormat = "#,##0"
$ws.Cells.Item(16, 8).Value = -3.030303030303
$ws.Cells.Item(16, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(16, 9).Value = 402
$ws.Cells.Item(16, 9).NumberFormat = "#,##0"
$ws.Cells.Item(16, 10).Value = 414
$ws.Cells.Item(16, 10).NumberFormat = "#,##0"
$ws.Cells.Item(16, 11).Value = -2.898550724637
$ws.Cells.Item(16, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(16, 12).Value = 25.233644859813
$ws.Cells.Item(16, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(16, 13).Value = 31.372549019607
$ws.Cells.Item(16, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(16, 14).Value = -67.475728155339
$ws.Cells.Item(16, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(17, 3).Value = 15
$ws.Cells.Item(17, 3).NumberFormat = "#,##0"
$ws.Cells.Item(17, 4).Value = 15
$ws.Cells.Item(17, 4).NumberFormat = "#,##0"
$ws.Cells.Item(17, 6).Value = 60
$ws.Cells.Item(17, 6).NumberFormat = "#,##0"
$ws.Cells.Item(17, 7).Value = 60
$ws.Cells.Item(17, 7).NumberFormat = "#,##0"
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(17, 9).Value = 667
$ws.Cells.Item(17, 9).NumberFormat = "#,##0"
$ws.Cells.Item(17, 10).Value = 645
$ws.Cells.Item(17, 10).NumberFormat = "#,##0"
$ws.Cells.Item(17, 11).Value = 3.410852713178
$ws.Cells.Item(17, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(17, 12).Value = 11.725293132328
$ws.Cells.Item(17, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(17, 13).Value = 117.97385620915
$ws.Cells.Item(17, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(17, 14).Value = -27.969762419006
$ws.Cells.Item(17, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(18, 3).Value = 3
$ws.Cells.Item(18, 3).NumberFormat = "#,##0"
$ws.Cells.Item(18, 4).Value = 2
$ws.Cells.Item(18, 4).NumberFormat = "#,##0"
$ws.Cells.Item(18, 5).Value = 50
$ws.Cells.Item(18, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(18, 6).Value = 16
$ws.Cells.Item(18, 6).NumberFormat = "#,##0"
$ws.Cells.Item(18, 8).Value = 14.285714285714
$ws.Cells.Item(18, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(18, 9).Value = 193
$ws.Cells.Item(18, 9).NumberFormat = "#,##0"
$ws.Cells.Item(18, 10).Value = 296
$ws.Cells.Item(18, 10).NumberFormat = "#,##0"
$ws.Cells.Item(18, 11).Value = -34.797297297297
$ws.Cells.Item(18, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(18, 12).Value = -1.530612244897
$ws.Cells.Item(18, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(18, 13).Value = 32.191780821917
$ws.Cells.Item(18, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(18, 14).Value = -81.70616113744
$ws.Cells.Item(18, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(19, 3).Value = 12
$ws.Cells.Item(19, 3).NumberFormat = "#,##0"
$ws.Cells.Item(19, 4).Value = 9
$ws.Cells.Item(19, 4).NumberFormat = "#,##0"
$ws.Cells.Item(19, 5).Value = 33.333333333333
$ws.Cells.Item(19, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(19, 6).Value = 42
$ws.Cells.Item(19, 6).NumberFormat = "#,##0"
$ws.Cells.Item(19, 7).Value = 35
$ws.Cells.Item(19, 7).NumberFormat = "#,##0"
$ws.Cells.Item(19, 8).Value = 20
$ws.Cells.Item(19, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(19, 9).Value = 490
$ws.Cells.Item(19, 9).NumberFormat = "#,##0"
$ws.Cells.Item(19, 10).Value = 488
$ws.Cells.Item(19, 10).NumberFormat = "#,##0"
$ws.Cells.Item(19, 11).Value = 0.409836065573
$ws.Cells.Item(19, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(19, 12).Value = 8.167770419426
$ws.Cells.Item(19, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(19, 13).Value = 117.777777777778
$ws.Cells.Item(19, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(19, 14).Value = 41.210374639769
$ws.Cells.Item(19, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(20, 3).Value = 8
$ws.Cells.Item(20, 3).NumberFormat = "#,##0"
$ws.Cells.Item(20, 4).Value = 4
$ws.Cells.Item(20, 4).NumberFormat = "#,##0"
$ws.Cells.Item(20, 5).Value = 100
$ws.Cells.Item(20, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(20, 6).Value = 16
$ws.Cells.Item(20, 6).NumberFormat = "#,##0"
$ws.Cells.Item(20, 8).Value = 6.666666666666
$ws.Cells.Item(20, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(20, 9).Value = 403
$ws.Cells.Item(20, 9).NumberFormat = "#,##0"
$ws.Cells.Item(20, 10).Value = 275
$ws.Cells.Item(20, 10).NumberFormat = "#,##0"
$ws.Cells.Item(20, 11).Value = 46.545454545454
$ws.Cells.Item(20, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(20, 12).Value = 84.018264840182
$ws.Cells.Item(20, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(20, 13).Value = 263.063063063063
$ws.Cells.Item(20, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(20, 14).Value = -10.840707964601
$ws.Cells.Item(20, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(21, 3).Value = 45
$ws.Cells.Item(21, 3).NumberFormat = "#,##0"
$ws.Cells.Item(21, 4).Value = 41
$ws.Cells.Item(21, 4).NumberFormat = "#,##0"
$ws.Cells.Item(21, 5).Value = 9.756097560975
$ws.Cells.Item(21, 5).NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Cells.Item(21, 6).Value = 167
$ws.Cells.Item(21, 6).NumberFormat = "#,##0"
$ws.Cells.Item(21, 7).Value = 162
$ws.Cells.Item(21, 7).NumberFormat = "#,##0"
$ws.Cells.Item(21, 8).Value = 3.086419753086
$ws.Cells.Item(21, 8).NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Cells.Item(21, 9).Value = 2198
$ws.Cells.Item(21, 9).NumberFormat = "#,##0"
$ws.Cells.Item(21, 10).Value = 2165
$ws.Cells.Item(21, 10).NumberFormat = "#,##0"
$ws.Cells.Item(21, 11).Value = 1.524249422632
$ws.Cells.Item(21, 11).NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Cells.Item(21, 12).Value = 19.521479064709
$ws.Cells.Item(21, 12).NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Cells.Item(21, 13).Value = 95.031055900621
$ws.Cells.Item(21, 13).NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Cells.Item(21, 14).Value = -46.663431206018
$ws.Cells.Item(21, 14).NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Cells.Item(23, 3).Value = 8
$ws.Cells.Item(23, 3).NumberFormat = "#,##0"
$ws.Cells.Item(23, 4).Value = 9
$ws.Cells.Item(23, 4).NumberFormat = "#,##0"
$ws.Cells.Item(23, 5).Value = -11.111111111111
$ws.Cells.Item(23, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(23, 9).Value = 397
$ws.Cells.Item(23, 9).NumberFormat = "#,##0"
$ws.Cells.Item(23, 10).Value = 348
$ws.Cells.Item(23, 10).NumberFormat = "#,##0"
$ws.Cells.Item(23, 11).Value = 14.080459770114
$ws.Cells.Item(23, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(23, 12).Value = 89.952153110047
$ws.Cells.Item(23, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(23, 13).Value = 108.947368421053
$ws.Cells.Item(23, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(24, 3).Value = 28
$ws.Cells.Item(24, 3).NumberFormat = "#,##0"
$ws.Cells.Item(24, 4).Value = 27
$ws.Cells.Item(24, 4).NumberFormat = "#,##0"
$ws.Cells.Item(24, 5).Value = 3.703703703703
$ws.Cells.Item(24, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(24, 6).Value = 86
$ws.Cells.Item(24, 6).NumberFormat = "#,##0"
$ws.Cells.Item(24, 8).Value = -5.494505494505
$ws.Cells.Item(24, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(24, 9).Value = 1124
$ws.Cells.Item(24, 9).NumberFormat = "#,##0"
$ws.Cells.Item(24, 10).Value = 1219
$ws.Cells.Item(24, 10).NumberFormat = "#,##0"
$ws.Cells.Item(24, 11).Value = -7.793273174733
$ws.Cells.Item(24, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(24, 12).Value = 24.198895027624
$ws.Cells.Item(24, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(24, 13).Value = 52.097428958051
$ws.Cells.Item(24, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(25, 3).Value = 14
$ws.Cells.Item(25, 3).NumberFormat = "#,##0"
$ws.Cells.Item(25, 4).Value = 13
$ws.Cells.Item(25, 4).NumberFormat = "#,##0"
$ws.Cells.Item(25, 5).Value = 7.692307692307
$ws.Cells.Item(25, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(25, 6).Value = 67
$ws.Cells.Item(25, 6).NumberFormat = "#,##0"
$ws.Cells.Item(25, 7).Value = 68
$ws.Cells.Item(25, 7).NumberFormat = "#,##0"
$ws.Cells.Item(25, 8).Value = -1.470588235294
$ws.Cells.Item(25, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(25, 9).Value = 987
$ws.Cells.Item(25, 9).NumberFormat = "#,##0"
$ws.Cells.Item(25, 10).Value = 924
$ws.Cells.Item(25, 10).NumberFormat = "#,##0"
$ws.Cells.Item(25, 11).Value = 6.818181818181
$ws.Cells.Item(25, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(25, 12).Value = 16.391509433962
$ws.Cells.Item(25, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(25, 13).Value = 17.081850533807
$ws.Cells.Item(25, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(26, 4).Value = 3
$ws.Cells.Item(26, 4).NumberFormat = "#,##0"
$ws.Cells.Item(26, 6).Value = 2
$ws.Cells.Item(26, 6).NumberFormat = "#,##0"
$ws.Cells.Item(26, 7).Value = 6
$ws.Cells.Item(26, 7).NumberFormat = "#,##0"
$ws.Cells.Item(26, 8).Value = -66.666666666666
$ws.Cells.Item(26, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(26, 10).Value = 62
$ws.Cells.Item(26, 10).NumberFormat = "#,##0"
$ws.Cells.Item(26, 11).Value = -17.741935483871
$ws.Cells.Item(26, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(27, 3).Value = 3
$ws.Cells.Item(27, 3).NumberFormat = "#,##0"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "0"
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = "***.*"
$ws.Cells.Item(27, 5).NumberFormat = "General"
$ws.Cells.Item(27, 6).Value = 6
$ws.Cells.Item(27, 6).NumberFormat = "#,##0"
$ws.Cells.Item(27, 7).Value = 3
$ws.Cells.Item(27, 7).NumberFormat = "#,##0"
$ws.Cells.Item(27, 8).Value = 100
$ws.Cells.Item(27, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(27, 9).Value = 92
$ws.Cells.Item(27, 9).NumberFormat = "#,##0"
$ws.Cells.Item(27, 11).Value = 27.777777777777
$ws.Cells.Item(27, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(27, 12).Value = 48.387096774193
$ws.Cells.Item(27, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "0"
$ws.Cells.Item(28, 3).NumberFormat = "General"
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(28, 4).NumberFormat = "#,##0"
$ws.Cells.Item(28, 5).Value = -100
$ws.Cells.Item(28, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(28, 7).Value = 4
$ws.Cells.Item(28, 7).NumberFormat = "#,##0"
$ws.Cells.Item(28, 8).Value = -75
$ws.Cells.Item(28, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(28, 10).Value = 36
$ws.Cells.Item(28, 10).NumberFormat = "#,##0"
$ws.Cells.Item(28, 11).Value = 5.555555555555
$ws.Cells.Item(28, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(28, 12).Value = -47.222222222222
$ws.Cells.Item(28, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(28, 13).Value = -19.148936170212
$ws.Cells.Item(28, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(28, 14).Value = -67.241379310344
$ws.Cells.Item(28, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(29, 3).NumberFormat = "@"
$ws.Cells.Item(29, 3).Value = "0"
$ws.Cells.Item(29, 3).NumberFormat = "General"
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(29, 4).NumberFormat = "#,##0"
$ws.Cells.Item(29, 5).Value = -100
$ws.Cells.Item(29, 5).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(29, 7).Value = 4
$ws.Cells.Item(29, 7).NumberFormat = "#,##0"
$ws.Cells.Item(29, 8).Value = -75
$ws.Cells.Item(29, 8).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(29, 10).Value = 32
$ws.Cells.Item(29, 10).NumberFormat = "#,##0"
$ws.Cells.Item(29, 11).Value = 0
$ws.Cells.Item(29, 11).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(29, 12).Value = -48.387096774193
$ws.Cells.Item(29, 12).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(29, 13).Value = -21.951219512195
$ws.Cells.Item(29, 13).NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Cells.Item(29, 14).Value = -71.171171171171
$ws.Cells.Item(29, 14).NumberFormat = "#,##0.0;""-""#,##0.0"
